$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Slightly narrow the sheet columns (pretest formatting tweak).
$ws.Columns("A").ColumnWidth = 5.166666666666667
$ws.Columns("B").ColumnWidth = 6.5
$ws.Columns("C").ColumnWidth = 9.333333333333334
$ws.Columns("D").ColumnWidth = 10.5
$ws.Columns("E").ColumnWidth = 9.166666666666666
$ws.Columns("F:G").ColumnWidth = 6.333333333333333
$ws.Columns("H").ColumnWidth = 8.0
$ws.Columns("I").ColumnWidth = 16.5
$ws.Columns("J").ColumnWidth = 5.333333333333333
$ws.Columns("K").ColumnWidth = 6.666666666666667
$ws.Columns("L:M").ColumnWidth = 7.333333333333333
$ws.Columns("N").ColumnWidth = 8.333333333333334
$ws.Columns("P").ColumnWidth = 13.0
$ws.Columns("Q").ColumnWidth = 8.5
$ws.Columns("R:S").ColumnWidth = 7.5
$ws.Columns("T").ColumnWidth = 6.166666666666667
$ws.Columns("U:V").ColumnWidth = 6.666666666666667
$ws.Columns("W").ColumnWidth = 8.333333333333334
$ws.Columns("X").ColumnWidth = 8.0
$ws.Columns("Y:Z").ColumnWidth = 8.333333333333334
$ws.Columns("AA").ColumnWidth = 27.333333333333332

# Update antenna count value in K4.
$ws.Range("K4").Value = 2
